# align opti with new reference
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the last data row (old row 10: sample 440 / 40h_NoOx_1 / S9) -
#    the reference set now stops at 40h_1 (row 9).
$ws.Rows.Item(10).Delete()

# 2) Update per-row numeric results (A/C/D/E/H) for the remaining 8 rows,
#    and flip the unique_anno boolean (F) + its numeric mirror (K) from
#    True/1 to False/0.
$rows = @(
    @{ Row=2; A=1512; C=99.27053892204847; D=74.027655903777;  E=4.440943291238329 },
    @{ Row=3; A=1513; C=99.31815534589504; D=74.30138568129331; E=4.457906322170901 },
    @{ Row=4; A=1514; C=98.68006956679432; D=72.13773775613342; E=5.911870705585674 },
    @{ Row=5; A=1515; C=99.29166895405103; D=74.22236228731563; E=4.060808259838651 },
    @{ Row=6; A=1516; C=98.82383710792956; D=74.13894165642355; E=4.857488448320816 },
    @{ Row=7; A=1517; C=98.79549216381972; D=74.13284942606347; E=4.076637407157326 },
    @{ Row=8; A=1518; C=99.09193504566448; D=73.9528143138853;  E=6.724284162620537 },
    @{ Row=9; A=1519; C=98.56170310665856; D=71.608583035256;   E=6.21327083988175 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $false
    $ws.Range("H$row").Value = 0.68
    $ws.Range("K$row").Value = 0
}

# 3) unique_anno_str (J2:J9) goes from the text "True" to "False". A plain
#    Value assignment of the word "False" auto-types to a boolean in this
#    engine (same literal-entry inference Excel itself does), so route the
#    text through a formula + paste-values round trip to land it as a real
#    string cell instead.
$jRange = $ws.Range("J2:J9")
$jRange.Formula = '="False"'
$jRange.Copy()
$jRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

Write-Host "done"
